# Apply weekly data-entry updates for three weekly sheets, then update the
# active-sheet/selection UI state to match the saved session.
#
# Sheets touched:
#   "24.10. - 30.10."  (week 6)
#   "31.10. - 6.11."   (week 7)
#   "7.11. - 13.11."   (week 8)
# The "Souhrn" sheet only contains SUM()/cross-sheet formulas, so it
# recalculates automatically once the underlying weekly values change.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "24.10. - 30.10."
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("24.10. - 30.10.")
$ws6.Range("C4").Value = 3
$ws6.Range("D4").Value = 1
$ws6.Range("G4").Value = 3
$ws6.Range("C6").Value = 2
$ws6.Range("D6").Value = 1
$ws6.Range("G6").Value = 1

# ---------------------------------------------------------------------
# 2) "31.10. - 6.11."
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("31.10. - 6.11.")

$ws7.Range("E3").Value = 0
$ws7.Range("F3").Value = 2
$ws7.Range("G3").Value = 3
$ws7.Range("H3").Value = 1
$ws7.Range("I3").Value = 7

$ws7.Range("C4").Value = 2
$ws7.Range("D4").Value = 5
$ws7.Range("E4").Value = 0
$ws7.Range("F4").Value = 0
$ws7.Range("G4").Value = 4
$ws7.Range("H4").Value = 1
$ws7.Range("I4").Value = 5

$ws7.Range("E5").Value = 5
$ws7.Range("F5").Value = 5
$ws7.Range("G5").Value = 5

$ws7.Range("C7").Value = 2
$ws7.Range("D7").Value = 4
$ws7.Range("E7").Value = 1
$ws7.Range("F7").Value = 4
$ws7.Range("G7").Value = 3
$ws7.Range("H7").Value = 1
$ws7.Range("I7").Value = 3

$ws7.Range("C8").Value = 2
$ws7.Range("D8").Value = 3
$ws7.Range("E8").Value = 2
$ws7.Range("F8").Value = 3
$ws7.Range("G8").Value = 3
$ws7.Range("H8").Value = 1
$ws7.Range("I8").Value = 3

# ---------------------------------------------------------------------
# 3) "7.11. - 13.11."
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("7.11. - 13.11.")

$ws8.Range("E3").Value = 1
$ws8.Range("G3").Value = 1
$ws8.Range("I3").Value = 1

$ws8.Range("C5").Value = 1
$ws8.Range("D5").Value = 4
$ws8.Range("E5").Value = 5
$ws8.Range("F5").Value = 1
$ws8.Range("G5").Value = 5
$ws8.Range("H5").Value = 3
$ws8.Range("I5").Value = 7

# ---------------------------------------------------------------------
# 4) Restore per-sheet selections (matches saved session state) and make
#    "7.11. - 13.11." the active tab/sheet last so it stays selected.
# ---------------------------------------------------------------------
$ws6.Activate()
$ws6.Range("H5").Select()

$ws7.Activate()
$ws7.Range("I21").Select()

$ws8.Activate()
$ws8.Range("I7").Select()
